$d = $word.ActiveDocument

# Remove the stray space between "${semester" and "}" so the merge field
# reads "${tipe} ${semester}" instead of "${tipe} ${semester }".
$d.Content.Find.Execute('${semester }', $true, $false, $false, $false, $false,
                         $true, 1, $false, '${semester}', 2)
